$wb = $excel.ActiveWorkbook

# --- ALC row 17 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 626.8246
$ws.Range("J17").Value = 626.8246
$ws.Range("L17").Value = 1880.4738
$ws.Range("N17").Value = -2216.4738

# --- ALC row 33 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 111863.22
$ws.Range("I33").Value = 143538.58
$ws.Range("J33").Value = 999.5
$ws.Range("K33").Value = 143538.58
$ws.Range("L33").Value = 999.5
$ws.Range("M33").Value = -143309.58
$ws.Range("N33").Value = -1457.5

# --- ALC row 132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 324699
$ws.Range("I132").Value = 2051.5
$ws.Range("J132").Value = 1430919
$ws.Range("K132").Value = 6154.5
$ws.Range("L132").Value = 4292757
$ws.Range("M132").Value = -3624.5
$ws.Range("N132").Value = -4297817

# --- ALC row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2684.25
$ws.Range("I138").Value = 754.0909
$ws.Range("J138").Value = 3228.6538
$ws.Range("K138").Value = 2262.2727
$ws.Range("L138").Value = 9685.9614
$ws.Range("M138").Value = 2877.7273
$ws.Range("N138").Value = -19965.9614

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2731087.2
$ws.Range("I32").Value = 7652.4165
$ws.Range("J32").Value = 12535453
$ws.Range("K32").Value = 7652.4165
$ws.Range("L32").Value = 12535453
$ws.Range("M32").Value = -7365.4165
$ws.Range("N32").Value = -12536027

# --- ARM row 74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1069.6818
$ws.Range("I74").Value = 1025.2667
$ws.Range("J74").Value = 1164.8572
$ws.Range("K74").Value = 1025.2667
$ws.Range("L74").Value = 1164.8572
$ws.Range("M74").Value = -151.2666999999999
$ws.Range("N74").Value = -2912.8572

# --- ARM row 77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1069.6818
$ws.Range("I77").Value = 1025.2667
$ws.Range("J77").Value = 1164.8572
$ws.Range("K77").Value = 5126.3335
$ws.Range("L77").Value = 5824.286
$ws.Range("M77").Value = -758.3334999999997
$ws.Range("N77").Value = -14560.286

# --- BSM row 134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 158363.69
$ws.Range("I134").Value = 5475
$ws.Range("J134").Value = 402985.6
$ws.Range("K134").Value = 16425
$ws.Range("L134").Value = 1208956.8
$ws.Range("M134").Value = -13890
$ws.Range("N134").Value = -1214026.8

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2987.4138
$ws.Range("I31").Value = 2776.4583
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 2776.4583
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -2481.4583
$ws.Range("N31").Value = -4590

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2987.4138
$ws.Range("I34").Value = 2776.4583
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 2776.4583
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -2574.4583
$ws.Range("N34").Value = -4404

# --- CRP row 134 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6443
$ws.Range("I134").Value = 5075.25
$ws.Range("J134").Value = 8266.666999999999
$ws.Range("K134").Value = 15225.75
$ws.Range("L134").Value = 24800.001
$ws.Range("M134").Value = -12690.75
$ws.Range("N134").Value = -29870.001

# --- CUL row 64 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 201142.2
$ws.Range("I64").Value = 906
$ws.Range("J64").Value = 334633
$ws.Range("K64").Value = 2718
$ws.Range("L64").Value = 1003899
$ws.Range("M64").Value = -2448
$ws.Range("N64").Value = -1004439

# --- CUL row 67 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 201142.2
$ws.Range("I67").Value = 906
$ws.Range("J67").Value = 334633
$ws.Range("K67").Value = 2718
$ws.Range("L67").Value = 1003899
$ws.Range("M67").Value = -1782
$ws.Range("N67").Value = -1005771

# --- CUL row 70 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 812
$ws.Range("I70").Value = 812
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2436
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2121
$ws.Range("N70").ClearContents()

# --- CUL row 73 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 812
$ws.Range("I73").Value = 812
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2436
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1344
$ws.Range("N73").ClearContents()

# --- CUL row 88 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 15500
$ws.Range("J88").Value = 28000
$ws.Range("L88").Value = 84000
$ws.Range("N88").Value = -84856

# --- CUL row 91 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 15500
$ws.Range("J91").Value = 28000
$ws.Range("L91").Value = 84000
$ws.Range("N91").Value = -86964

# --- CUL row 94 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 5266.6665
$ws.Range("I94").Value = 3400
$ws.Range("J94").Value = 9000
$ws.Range("K94").Value = 10200
$ws.Range("L94").Value = 27000
$ws.Range("M94").Value = -9524
$ws.Range("N94").Value = -28352

# --- CUL row 95 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 7400

# --- CUL row 97 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 23810322
$ws.Range("I97").Value = 35714484
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 107143452
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -107142956
$ws.Range("N97").Value = -6992

# --- CUL row 126 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 4443.75
$ws.Range("I126").Value = 1550
$ws.Range("J126").Value = 5408.3335
$ws.Range("K126").Value = 4650
$ws.Range("L126").Value = 16225.0005
$ws.Range("M126").Value = 290
$ws.Range("N126").Value = -26105.0005

# --- GSM row 62 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# --- GSM row 65 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# --- GSM row 74 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 49465.5
$ws.Range("J74").Value = 49465.5
$ws.Range("L74").Value = 49465.5
$ws.Range("N74").Value = -51337.5

# --- GSM row 77 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 49465.5
$ws.Range("J77").Value = 49465.5
$ws.Range("L77").Value = 148396.5
$ws.Range("N77").Value = -157756.5

# --- GSM row 86 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# --- GSM row 88 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 65000
$ws.Range("J88").Value = 65000
$ws.Range("L88").Value = 65000
$ws.Range("N88").Value = -65902

# --- GSM row 89 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# --- GSM row 91 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H91").Value = 65000
$ws.Range("J91").Value = 65000
$ws.Range("L91").Value = 65000
$ws.Range("N91").Value = -68120
